$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-05-28 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-05-29 Thursday", 2)

# The single table holds the division problems in 5 "data" rows
# (1, 5, 9, 13, 17) with 5 columns each, for 25 problems total,
# followed by blank rows used for writing the answer.
$t = $d.Tables.Item(1)

$replacements = @(
    "31÷3=", "20÷3=",
    "50÷8=", "11÷3=",
    "13÷2=", "59÷2=",
    "21÷2=", "19÷6=",
    "18÷8=", "77÷2=",

    "27÷8=", "84÷3=",
    "89÷7=", "74÷2=",
    "21÷6=", "46÷2=",
    "49÷5=", "11÷7=",
    "69÷5=", "42÷8=",

    "62÷6=", "81÷3=",
    "26÷8=", "23÷3=",
    "99÷4=", "96÷6=",
    "83÷6=", "26÷9=",
    "72÷6=", "10÷3=",

    "53÷7=", "72÷6=",
    "89÷2=", "31÷3=",
    "13÷2=", "82÷9=",
    "67÷9=", "57÷8=",
    "44÷3=", "43÷8=",

    "34÷2=", "87÷3=",
    "82÷7=", "23÷4=",
    "18÷8=", "43÷9=",
    "22÷2=", "89÷2=",
    "22÷4=", "71÷9="
)

$dataRows = @(1, 5, 9, 13, 17)
$idx = 0
foreach ($row in $dataRows) {
    for ($col = 1; $col -le 5; $col++) {
        $oldVal = $replacements[$idx]
        $newVal = $replacements[$idx + 1]
        $idx = $idx + 2

        $cell = $t.Cell($row, $col)
        $r = $cell.Range
        $r.End = $r.End - 1
        $r.Find.Execute($oldVal, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newVal, 2)
    }
}
